$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$newSheet.Name = "TestLogin"
$newSheet.Range("A1").Value = "Result"
$newSheet.Range("A1").NumberFormat = "@"
$newSheet.Range("A1").WrapText = $true
$newSheet.Range("A1").Borders.LineStyle = 1
